# Add ability to get data and send in sms: new "subject_name" text prompt
# is inserted as row 2 of the "survey" sheet, and the "survey" tab becomes
# the active tab/selection (instead of "settings").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new row at position 2; this pushes the existing row 2
# (send_sms / send_sms / ...) down to row 3.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the subject_name text prompt.
# Set column B before column A so new shared strings are appended to the
# table in the same order as the target workbook (subject_name, text, ...).
$ws.Cells.Item(2,2).Value = "subject_name"
$ws.Cells.Item(2,1).Value = "text"
$ws.Cells.Item(2,3).Value = "Enter the subject's name."

# Match the style used by the header row's "display.text" cell (C1).
$styleSrc = $ws.Cells.Item(1,3)
$ws.Cells.Item(2,1).Style = $styleSrc.Style
$ws.Cells.Item(2,2).Style = $styleSrc.Style
$ws.Cells.Item(2,3).Style = $styleSrc.Style

$ws.Rows.Item(2).RowHeight = 12

# Make the "survey" sheet the active tab, with D3 selected (instead of the
# previously active "settings" sheet).
$ws.Activate() | Out-Null
$ws.Range("D3").Select() | Out-Null
